$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "d_thomson" in F1, matching the style of the other headers (E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "d_thomson"

# Updated sigma (C), d_sigma (D) values, refreshed thomson (E) values,
# and new d_thomson (F) values for rows 2-8
$ws.Cells.Item(2, 3).Value = 0.000000000000000000000000000002706292830616987
$ws.Cells.Item(2, 4).Value = 0.00000000000000000000000000000004169523577880892
$ws.Cells.Item(2, 5).Value = 0.00000000000000000000000000006652458721354302
$ws.Cells.Item(2, 6).Value = 0.00000000000000000000000000000000000004089603650513666

$ws.Cells.Item(3, 3).Value = 0.000000000000000000000000000001671436234244726
$ws.Cells.Item(3, 4).Value = 0.00000000000000000000000000000005659619734661688
$ws.Cells.Item(3, 5).Value = 0.00000000000000000000000000006652458721354302
$ws.Cells.Item(3, 6).Value = 0.00000000000000000000000000000000000004089603650513666

$ws.Cells.Item(4, 3).Value = 0.0000000000000000000000000000005447411529868519
$ws.Cells.Item(4, 4).Value = 0.00000000000000000000000000000001324136866993893
$ws.Cells.Item(4, 5).Value = 0.00000000000000000000000000006652458721354302
$ws.Cells.Item(4, 6).Value = 0.00000000000000000000000000000000000004089603650513666

$ws.Cells.Item(5, 3).Value = 0.0000000000000000000000000000002997252158003094
$ws.Cells.Item(5, 4).Value = 0.00000000000000000000000000000001296765517649
$ws.Cells.Item(5, 5).Value = 0.00000000000000000000000000006652458721354302
$ws.Cells.Item(5, 6).Value = 0.00000000000000000000000000000000000004089603650513666

$ws.Cells.Item(6, 3).Value = 0.0000000000000000000000000000002490109348451165
$ws.Cells.Item(6, 4).Value = 0.00000000000000000000000000000001609715051891146
$ws.Cells.Item(6, 5).Value = 0.00000000000000000000000000006652458721354302
$ws.Cells.Item(6, 6).Value = 0.00000000000000000000000000000000000004089603650513666

$ws.Cells.Item(7, 3).Value = 0.0000000000000000000000000000002467164843893904
$ws.Cells.Item(7, 4).Value = 0.00000000000000000000000000000002130347107798132
$ws.Cells.Item(7, 5).Value = 0.00000000000000000000000000006652458721354302
$ws.Cells.Item(7, 6).Value = 0.00000000000000000000000000000000000004089603650513666

$ws.Cells.Item(8, 3).Value = 0.0000000000000000000000000000001799017339060681
$ws.Cells.Item(8, 4).Value = 0.00000000000000000000000000000004081824766623204
$ws.Cells.Item(8, 5).Value = 0.00000000000000000000000000006652458721354302
$ws.Cells.Item(8, 6).Value = 0.00000000000000000000000000000000000004089603650513666

